$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Metal Door asset modeled: record its poly count in the "Poly Count" column (D4)
$ws.Range("D4").Value = "64 Polygons"
